# Apply updates to MMFAD93_test_1.xlsx per commit "update bugs in MMFAD9A"
$wb = $excel.ActiveWorkbook

# --- Sheet 1: question_answers (text answers in column B) ---
$ws1 = $wb.Worksheets.Item("question_answers")

$qaUpdates = @{
    4  = "5"
    6  = "4"
    7  = "1"
    8  = "4"
    10 = "5"
    12 = "4"
    14 = "3"
    15 = "5"
    16 = "2"
    17 = "4"
    18 = "1"
    19 = "4"
    20 = "3"
    21 = "2"
    23 = "4"
    24 = "4"
    26 = "4"
    27 = "4"
    28 = "3"
    29 = "5"
    30 = "2"
    31 = "4"
    32 = "4"
    33 = "2"
    35 = "5"
    36 = "2"
    37 = "1"
    38 = "3"
    39 = "1"
    41 = "5"
    42 = "1"
    43 = "5"
    45 = "3"
    46 = "3"
    47 = "2"
    49 = "4"
    50 = "5"
    51 = "1"
    52 = "4"
    54 = "3"
}

foreach ($row in $qaUpdates.Keys) {
    $cell = $ws1.Range("B$row")
    # Keep these as text values (original cells are stored as text/inline strings)
    $cell.NumberFormat = "@"
    $cell.Value = $qaUpdates[$row]
}

# --- Sheet 2: outputs (numeric scores in column B) ---
$ws2 = $wb.Worksheets.Item("outputs")

$outputsUpdates = @{
    2 = 8
    3 = 26
    4 = 17
    5 = 18
    6 = 21
    7 = 20
    9 = 143
}

foreach ($row in $outputsUpdates.Keys) {
    $ws2.Range("B$row").Value = $outputsUpdates[$row]
}
